# #67 Ajouter l'année 2020 dans la fiche de stats
#
# The template shifts its rolling 4-season window forward by one year:
#   2019 -> 2020 (cols D:F), 2018 -> 2019 (cols G:I),
#   2017 -> 2018 (cols J:L), 2016 -> 2017 (cols M:O)
# on both "Par saison" worksheets, and the "Général" sheet's
# "Version du modèle" date moves forward too.

$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item("Par saison (fin de saison)"), $wb.Worksheets.Item("Par saison (date de génération)"))) {
    # Row 3: season range headers (merged cells, value lives in the first cell)
    $ws.Range("D3").Value = "2020-2019"
    $ws.Range("G3").Value = "2019-2018"
    $ws.Range("J3").Value = "2018-2019"
    $ws.Range("M3").Value = "2017-2018"

    # Row 5: JSP template placeholders for each season's "responsables"/"jeunes"
    $ws.Range("D5").Value = '${effectif.responsables.2020}'
    $ws.Range("E5").Value = '${effectif.jeunes.2020}'
    $ws.Range("G5").Value = '${effectif.responsables.2019}'
    $ws.Range("H5").Value = '${effectif.jeunes.2019}'
    $ws.Range("J5").Value = '${effectif.responsables.2018}'
    $ws.Range("K5").Value = '${effectif.jeunes.2018}'
    $ws.Range("M5").Value = '${effectif.responsables.2017}'
    $ws.Range("N5").Value = '${effectif.jeunes.2017}'
}

# "Général" sheet: bump the model version date forward (2019-11-22 -> 2020-09-06)
$wsGeneral = $wb.Worksheets.Item("Général")
$wsGeneral.Range("B1").Value = 44080
